# Update cryptos list values per upstream data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.488.79'
$ws.Range("E2").Value = '  -2.10%  '
$ws.Range("D3").Value = '2.422.59'
$ws.Range("E3").Value = '  -1.75%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = '''550.39'
$ws.Range("E5").Value = '  -1.34%  '
$ws.Range("D6").Value = '''158.38'
$ws.Range("E6").Value = '  -2.60%  '
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("D8").Value = '''0.505'
$ws.Range("E8").Value = '  +0.25%  '
$ws.Range("D9").Value = '''0.157'
$ws.Range("E9").Value = '  +4.57%  '
$ws.Range("E10").Value = '  -1.02%  '
$ws.Range("D11").Value = '''0.327'
$ws.Range("E11").Value = '  -1.97%  '
$ws.Range("D12").Value = '''4.76'
$ws.Range("E12").Value = '  -1.28%  '
$ws.Range("D13").Value = '67.430.69'
$ws.Range("E13").Value = '  -1.98%  '
$ws.Range("E14").Value = '  -0.18%  '
$ws.Range("D15").Value = '''22.88'
$ws.Range("E15").Value = '  -2.80%  '
$ws.Range("D16").Value = '''10.30'
$ws.Range("E16").Value = '  -4.15%  '
$ws.Range("D17").Value = '''327.62'
$ws.Range("E17").Value = '  -4.18%  '
$ws.Range("D18").Value = '''6.81'
$ws.Range("E18").Value = '  -3.23%  '
$ws.Range("D19").Value = '''3.76'
$ws.Range("E19").Value = '  -0.72%  '
$ws.Range("E20").Value = '  -0.40%  '
$ws.Range("E21").Value = '  -1.63%  '
$ws.Range("D22").Value = '''65.86'
$ws.Range("E22").Value = '  -1.51%  '
$ws.Range("D23").Value = '''3.59'
$ws.Range("E23").Value = '  -2.32%  '
$ws.Range("D24").Value = '''8.00'
$ws.Range("E24").Value = '  -1.57%  '
$ws.Range("D25").Value = '0.0₃0796'
$ws.Range("E25").Value = '  -2.83%  '
$ws.Range("D26").Value = '''6.99'
$ws.Range("E26").Value = '  -2.64%  '
$ws.Range("E27").Value = '  +0.03%  '
$ws.Range("D28").Value = '''413.12'
$ws.Range("E28").Value = '  -5.67%  '
$ws.Range("D29").Value = '''1.11'
$ws.Range("E29").Value = '  -2.16%  '
$ws.Range("D30").Value = '''1.58'
$ws.Range("E30").Value = '  -1.81%  '
$ws.Range("D31").Value = '''160.06'
$ws.Range("E31").Value = '  +1.71%  '
$ws.Range("D32").Value = '''18.92'
$ws.Range("E32").Value = '  -0.68%  '
$ws.Range("E33").Value = '  -0.09%  '
$ws.Range("D34").Value = '''17.73'
$ws.Range("E34").Value = '  -0.74%  '
$ws.Range("E35").Value = '  -4.66%  '
$ws.Range("D36").Value = '''0.292'
$ws.Range("E36").Value = '  -3.42%  '
$ws.Range("D37").Value = '''4.21'
$ws.Range("E37").Value = '  -5.17%  '
$ws.Range("D38").Value = '''1.44'
$ws.Range("E38").Value = '  -2.07%  '
$ws.Range("E39").Value = '  -3.44%  '
$ws.Range("B40").Value = 'Aave'
$ws.Range("C40").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D40").Value = '''129.56'
$ws.Range("E40").Value = '  -2.45%  '
$ws.Range("B41").Value = 'dogwifhat'
$ws.Range("C41").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D41").Value = '''1.97'
$ws.Range("E41").Value = '  -4.77%  '
$ws.Range("B42").Value = 'Filecoin'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D42").Value = '''3.29'
$ws.Range("E42").Value = '  -1.63%  '
$ws.Range("E43").Value = '  -1.22%  '
$ws.Range("D44").Value = '''0.474'
$ws.Range("E44").Value = '  -1.95%  '
$ws.Range("E45").Value = '  -1.91%  '
$ws.Range("D46").Value = '''0.0911'
$ws.Range("E46").Value = '  +0.34%  '
$ws.Range("E47").Value = '  +0.13%  '
$ws.Range("D48").Value = '''1.32'
$ws.Range("E48").Value = '  -8.19%  '
$ws.Range("D49").Value = '''16.38'
$ws.Range("E49").Value = '  -3.02%  '
$ws.Range("D50").Value = '0.0₆0201'
$ws.Range("E50").Value = '  -0.69%  '
$ws.Range("D51").Value = '''0.0425'
$ws.Range("E51").Value = '  -1.09%  '
